$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the existing data row 3 onto the
# previously-empty row 6 so the new test case visually matches the
# other rows (TestCase / Username / Password / Email columns).
$ws.Range("A3:D3").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)

# Fill in the new "repeated user registration" test case values.
$ws.Range("A6").Value = "TC_13"
$ws.Range("B6").Value = "KaranPrinja"
$ws.Range("C6").Value = "Test@12345"
$ws.Range("D6").Value = "karanprinja@hotmail.com"

# Style the e-mail address like the other e-mail cells (underlined,
# size 10, Helvetica Neue, "hyperlink" colored) by formatting the
# run of characters (done in two passes so the formatting is stored
# as explicit run-level rich text, matching the other rows).
$emailLen = 23
$part1 = $ws.Range("D6").Characters(1, $emailLen - 1)
$part1.Font.Underline = 2
$part1.Font.Size = 10
$part1.Font.Color = 16776960
$part1.Font.Name = "Helvetica Neue"

$part2 = $ws.Range("D6").Characters($emailLen, 1)
$part2.Font.Underline = 2
$part2.Font.Size = 10
$part2.Font.Color = 16776960
$part2.Font.Name = "Helvetica Neue"

# Add the mailto: hyperlink (Test Listener) for the new e-mail cell,
# just like the other TestCase rows.
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:karanprinja@hotmail.com", "", "", "karanprinja@hotmail.com")

# Adding the hyperlink re-stamps the cell with its own "Hyperlink"
# style; restore the column's normal data-row style (same as C6) so
# the cell formatting matches the rest of the table.
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)

# Clean up the now-unused built-in "Hyperlink" cell style that was
# auto-created by Hyperlinks.Add above.
$wb.Styles.Item("Hyperlink").Delete()
